$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.524.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -2.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3904"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08357"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.105"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.77"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.225"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.881.72"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.42"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.284"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.913"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.551.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.15"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.226"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.095.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.47"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.415"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.040"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.740"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02451"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06555"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.933"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2165"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.026"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.180"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.236"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6379"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6002"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.99"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.690"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.002"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.218"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.99"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.139"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.21%  "
